# Fruta / hortaliza, semanal
# Insert two new weekly observation rows before the current row 96
# (this pushes the former rows 96-99 down to rows 98-101, unchanged),
# then populate the two newly inserted rows (96, 97) with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("96:97").Insert()

# New row 96
$ws.Cells.Item(96, 1).Value = 4
$ws.Cells.Item(96, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(96, 3).Value = "Los Lagos"
$ws.Cells.Item(96, 4).Value = 44509
$ws.Cells.Item(96, 5).Value = 10
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100108
$ws.Cells.Item(96, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(96, 9).Value = 100108002
$ws.Cells.Item(96, 10).Value = "Mango"
$ws.Cells.Item(96, 11).Value = "Sin especificar"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 300
$ws.Cells.Item(96, 14).Value = 8000
$ws.Cells.Item(96, 15).Value = 8500
$ws.Cells.Item(96, 16).Value = 8250
$ws.Cells.Item(96, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(96, 18).Value = "Perú"
$ws.Cells.Item(96, 19).Value = 2062
$ws.Cells.Item(96, 20).Value = 4

# New row 97
$ws.Cells.Item(97, 1).Value = 4
$ws.Cells.Item(97, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(97, 3).Value = "Los Lagos"
$ws.Cells.Item(97, 4).Value = 44509
$ws.Cells.Item(97, 5).Value = 10
$ws.Cells.Item(97, 6).Value = "Fruta"
$ws.Cells.Item(97, 7).Value = 100108
$ws.Cells.Item(97, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(97, 9).Value = 100108002
$ws.Cells.Item(97, 10).Value = "Mango"
$ws.Cells.Item(97, 11).Value = "Sin especificar"
$ws.Cells.Item(97, 12).Value = "Segunda"
$ws.Cells.Item(97, 13).Value = 60
$ws.Cells.Item(97, 14).Value = 6000
$ws.Cells.Item(97, 15).Value = 6000
$ws.Cells.Item(97, 16).Value = 6000
$ws.Cells.Item(97, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(97, 18).Value = "Perú"
$ws.Cells.Item(97, 19).Value = 1500
$ws.Cells.Item(97, 20).Value = 4
